$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold + border + center/top alignment) from H1
# onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J (identical per row, rows 2-18)
$values = @(9, 8, 9, 8, 9, 8, 8, 8, 8, 8, 8, 8, 8, 6, 6, 4, 7)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $val = $values[$i]
    $ws.Cells.Item($row, 9).Value = $val
    $ws.Cells.Item($row, 10).Value = $val
}
